$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Week 3 task summary data ---
# New shared strings must be introduced in this order so the resulting
# sharedStrings table matches: Project Planning, Analysis/requirements
# Elicitation, Finalise Project Plan, Client Meeting..., Richard Dobson.
$ws.Range("A3").Value = "Project Planning"
$ws.Range("A4").Value = "Analysis/requirements Elicitation"
$ws.Range("B3").Value = "Finalise Project Plan"
$ws.Range("B4").Value = "Client Meeting, Gather and analyse requirements"
$ws.Range("C1").Value = "Richard Dobson"

$ws.Range("E1").Value = 3

$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 18
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2

$ws.Range("B14").Value = 60

# --- Widen the "Stage" column ---
$ws.Columns.Item(1).ColumnWidth = 31.1

# --- Update the active selection to C1 ---
$null = $ws.Range("C1").Select()
